# EcoSysEM package layout.docx - apply the two text edits from the commit.
#
# Both edits change a run boundary inside runs of code-style text: editing
# text through the Word object model automatically re-merges any two
# adjacent runs that end up with identical run formatting (exactly like a
# save-time "normalise runs" pass in Word itself). To land the intended
# run split / run merge precisely - instead of a coarser merge that would
# also sweep in neighbouring, differently-rsid'd runs that share the same
# formatting - each edit is done in two scoped steps: touch the piece that
# must stay put last, and briefly diverge the formatting of the other
# piece while the textual edit lands so it is not folded back together,
# then restore that formatting immediately after.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "plotDeltaGr"  ->  "export" + "DeltaGr"   (bold, unchanged rPr)
# ---------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("plotDeltaGr", $true, $false, $false, $false, $false, $true)
$fullStart = $rng.Start
$fullEnd = $rng.End

# Diverge the "DeltaGr" tail's colour so the head edit below won't cause it
# to be swept back into a single run.
$tail = $d.Range($fullStart + 4, $fullEnd)
$tail.Font.Color = 255

# Shrink "plot" -> "export" in place.
$head = $d.Range($fullStart, $fullStart + 4)
$head.Text = "export"

# Restore the tail's colour (fresh Range so the write is scoped to just the
# "DeltaGr" text, leaving "export" and everything before/after untouched).
$tail2 = $d.Range($head.End, $head.End + 7)
$tail2.Font.Color = 4210752

# ---------------------------------------------------------------------
# Edit 2: "ecosysem_" + "spyder" + ".py"  ->  "ecosysem_" + "spyder.py"
# (the "spyder" / ".py" runs merge into one run; "ecosysem_" is untouched)
# ---------------------------------------------------------------------
$rng2 = $d.Content
$null = $rng2.Find.Execute("ecosysem_spyder.py", $true, $false, $false, $false, $false, $true)
$fullStart2 = $rng2.Start
$fullEnd2 = $rng2.End

# Diverge the "ecosysem_" prefix's colour so the tail edit below won't pull
# it into the merge.
$prefix = $d.Range($fullStart2, $fullStart2 + 9)
$prefix.Font.Color = 255

# Force a genuine content diff over "spyder.py" (via a temp marker) so the
# "spyder" and ".py" runs actually collapse into a single run.
$tailName = $d.Range($fullStart2 + 9, $fullEnd2)
$tailName.Text = "spyder_py_TEMP_MARKER"
$tailName.Text = "spyder.py"

# Restore the prefix's colour (fresh Range scoped to just "ecosysem_").
$prefix2 = $d.Range($fullStart2, $fullStart2 + 9)
$prefix2.Font.Color = 4210752
